$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.265.55"
$ws.Range("E2").Value = "  -3.48%  "

$ws.Range("D3").Value = "3.182.89"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.07"
$ws.Range("E5").Value = "  +0.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.22"
$ws.Range("E6").Value = "  -6.42%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.151.67"
$ws.Range("E8").Value = "  -2.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  -3.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("E10").Value = "  -7.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("E11").Value = "  -4.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.473"
$ws.Range("E12").Value = "  -5.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  -6.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.94"
$ws.Range("E14").Value = "  -7.22%  "

$ws.Range("D15").Value = "3.663.87"
$ws.Range("E15").Value = "  -2.80%  "

$ws.Range("D16").Value = "64.187.04"
$ws.Range("E16").Value = "  -3.65%  "

$ws.Range("E17").Value = "  +1.20%  "

$ws.Range("D18").Value = "3.163.02"
$ws.Range("E18").Value = "  -2.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.93"
$ws.Range("E19").Value = "  -4.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.68"
$ws.Range("E20").Value = "  -5.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.55"
$ws.Range("E21").Value = "  -4.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.713"
$ws.Range("E22").Value = "  -3.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.72"
$ws.Range("E23").Value = "  -3.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.72"
$ws.Range("E24").Value = "  -5.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.49"
$ws.Range("E25").Value = "  -3.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.88"
$ws.Range("E27").Value = "  -3.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.38"
$ws.Range("E28").Value = "  -7.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -7.02%  "

$ws.Range("E30").Value = "  -28.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.77"
$ws.Range("E31").Value = "  -1.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.74"
$ws.Range("E32").Value = "  -5.60%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.15"
$ws.Range("E34").Value = "  -7.13%  "

$ws.Range("E35").Value = "  -5.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.00"
$ws.Range("E36").Value = "  -5.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.22"
$ws.Range("E37").Value = "  -1.99%  "

$ws.Range("D38").Value = "0.0₃0716"
$ws.Range("E38").Value = "  -9.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "451.39"
$ws.Range("E39").Value = "  -8.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.90"
$ws.Range("E40").Value = "  -7.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0395"
$ws.Range("E41").Value = "  -5.97%  "

$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.42"
$ws.Range("E42").Value = "  -3.11%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.118"
$ws.Range("E43").Value = "  -8.17%  "

$ws.Range("D44").Value = "2.851.21"
$ws.Range("E44").Value = "  -3.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.268"
$ws.Range("E45").Value = "  -8.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.24"
$ws.Range("E46").Value = "  -8.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.31"
$ws.Range("E47").Value = "  -6.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.114"
$ws.Range("E49").Value = "  -3.94%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.30"
$ws.Range("E50").Value = "  -4.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.58"
$ws.Range("E51").Value = "  -1.73%  "
